# Generate Report for Handback
# For the "7a966149-f594-49f7-841e-1201304368db" row (row 6) on both the
# "zh-cn" and "de-de" sheets, fill in the target-file / handback-file /
# handback-datetime / error-detail columns, add a hyperlink on the new
# "Latest Target File" cell that mirrors the existing handoff-file
# hyperlink in column A, and widen the affected columns to match.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/656dd222ca7c0cea26de5be5041b009b4de2783e/e2e/7a966149-f594-49f7-841e-1201304368db.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc5c06a254379b2c7a0420aa0f2901b33553054f/e2e/7a966149-f594-49f7-841e-1201304368db.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/656dd222ca7c0cea26de5be5041b009b4de2783e/e2e/7a966149-f594-49f7-841e-1201304368db.md."

function Update-LanguageSheet {
    param($ws, $handbackFileName, $handbackDateTime)

    # Latest Target File (I6) — new hyperlink pointing at the same handoff
    # markdown file as A6, displaying the .md file name.
    $ws.Hyperlinks.Add($ws.Range("I6"), $targetUrl, "", "", "7a966149-f594-49f7-841e-1201304368db.md") | Out-Null
    $ws.Range("I6").Style = $ws.Range("A6").Style

    # Latest Handback File (J6)
    $ws.Range("J6").Value = $handbackFileName

    # Latest Handback DateTime (K6)
    $ws.Range("K6").Value = $handbackDateTime

    # Error Detail (P6)
    $ws.Range("P6").Value = $errorDetail

    # Widen the columns that now hold longer content.
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
    $ws.Columns.Item(16).ColumnWidth = 40
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LanguageSheet $wsZhCn "7a966149-f594-49f7-841e-1201304368db.c31be0b884059c8b18c505418ae27ea9d0b7d375.zh-cn.xlf" "2016-08-31 11:54:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LanguageSheet $wsDeDe "7a966149-f594-49f7-841e-1201304368db.c31be0b884059c8b18c505418ae27ea9d0b7d375.de-de.xlf" "2016-08-31 11:54:38"
